# Update countries & provincias Spain
# - Update the "last updated" timestamp text
# - Update COVID stats for several countries; Peru overtakes India and
#   Sudafrica overtakes Egipto in the ranking, so those rows swap places
# - Refresh the numeric statistics for several other countries

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Title / last updated timestamp (row 1)
$ws.Range("A1").Value = "Datos actualizados a 6 de Mayo de 2020 a las 22:33"

# Row 4 - Estados Unidos: refreshed totals
$ws.Range("B4").Value = 1252689
$ws.Range("C4").Value = 15056
$ws.Range("D4").Value = 205215
$ws.Range("E4").Value = 973677
$ws.Range("G4").Value = 1526
$ws.Range("H4").Value = 73797

# Rows 16-17 - Peru overtakes India, rows swap
$ws.Range("A16").Value = "Peru"
$ws.Range("B16").Value = 54817
$ws.Range("C16").Value = 3628
$ws.Range("D16").Value = 17527
$ws.Range("E16").Value = 35757
$ws.Range("F16").Value = 717
$ws.Range("G16").Value = 89
$ws.Range("H16").Value = 1533

$ws.Range("A17").Value = "India"
$ws.Range("B17").Value = 52987
$ws.Range("C17").Value = 3587
$ws.Range("D17").Value = 15331
$ws.Range("E17").Value = 35871
$ws.Range("F17").Value = 0
$ws.Range("G17").Value = 92
$ws.Range("H17").Value = 1785

# Row 22 - Suiza: refreshed totals
$ws.Range("D22").Value = 25700
$ws.Range("E22").Value = 2555

# Row 32 - Israel: refreshed totals
$ws.Range("B32").Value = 16310
$ws.Range("C32").Value = 21
$ws.Range("D32").Value = 10637
$ws.Range("E32").Value = 5434
$ws.Range("F32").Value = 89
$ws.Range("G32").Value = 1
$ws.Range("H32").Value = 239

# Rows 49-50 - Sudafrica overtakes Egipto, rows swap
$ws.Range("A49").Value = "Sudafrica"
$ws.Range("B49").Value = 7808
$ws.Range("C49").Value = 236
$ws.Range("D49").Value = 3153
$ws.Range("E49").Value = 4502
$ws.Range("F49").Value = 36
$ws.Range("G49").Value = 5
$ws.Range("H49").Value = 153

$ws.Range("A50").Value = "Egipto"
$ws.Range("B50").Value = 7588
$ws.Range("C50").Value = 387
$ws.Range("D50").Value = 1815
$ws.Range("E50").Value = 5304
$ws.Range("F50").Value = 0
$ws.Range("G50").Value = 17
$ws.Range("H50").Value = 469

# Row 98 - Sudan: refreshed totals
$ws.Range("E98").Value = 723
$ws.Range("G98").Value = 4
$ws.Range("H98").Value = 49

# Row 104 - Costa Rica: refreshed totals
$ws.Range("B104").Value = 761
$ws.Range("C104").Value = 6
$ws.Range("D104").Value = 428
$ws.Range("E104").Value = 327
